$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Website Chals - B9 gets a right-aligned style
# ---------------------------------------------------------------------
$wsWebsite = $wb.Worksheets.Item("Website Chals")
$wsWebsite.Range("B9").HorizontalAlignment = -4152

# ---------------------------------------------------------------------
# Sheet: Wargame Info - remove the now-unused "Points" column (column E)
# which shifts the old F/G columns left into E/F, matching the diff.
# ---------------------------------------------------------------------
$wsWargame = $wb.Worksheets.Item("Wargame Info")
$wsWargame.Columns("E:E").Delete()

# Update the challenge counts for the coaches handout
$wsWargame.Range("B18").Value = 5
$wsWargame.Range("B21").Value = 7
$wsWargame.Range("B23").Value = 8

# New note about the coaches handout file (adds a new shared string)
$wsWargame.Range("C24").Value = "Waiting on Corey's…"

# ---------------------------------------------------------------------
# Sheet view / selection bookkeeping
# ---------------------------------------------------------------------

# Challenge Flags: scroll back to the top and select F8:G8
$wsFlags = $wb.Worksheets.Item("Challenge Flags")
$wsFlags.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$wsFlags.Range("F8:G8").Select()

# Website Chals: scroll so column C is left-most and select E27
$wsWebsite.Activate()
$excel.ActiveWindow.ScrollColumn = 3
$wsWebsite.Range("E27").Select()

# Scoreboard: keep its own selection untouched, just make sure it is no
# longer the active tab (handled automatically once another sheet is
# activated afterwards).
$wsScoreboard = $wb.Worksheets.Item("Scoreboard")
$wsScoreboard.Activate()
$wsScoreboard.Range("A3").Select()

# Wargame Info: becomes the active tab; scroll back to the top and select C7
$wsWargame.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$wsWargame.Range("C7").Select()
